$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sector breakdown values to reflect refreshed data pull
# Row 2
$ws.Cells.Item(2, "B").Value = 0
$ws.Cells.Item(2, "C").Value = 0
$ws.Cells.Item(2, "D").Value = 0
$ws.Cells.Item(2, "E").Value = 0
$ws.Cells.Item(2, "F").Value = 0
$ws.Cells.Item(2, "G").Value = 0
$ws.Cells.Item(2, "H").Value = 0
$ws.Cells.Item(2, "I").Value = 0
$ws.Cells.Item(2, "K").Value = 0
$ws.Cells.Item(2, "L").Value = 0
$ws.Cells.Item(2, "M").Value = 0

# Row 3
$ws.Cells.Item(3, "B").Value = 9.934
$ws.Cells.Item(3, "C").Value = 9.171
$ws.Cells.Item(3, "D").Value = 29.139
$ws.Cells.Item(3, "E").Value = 6.129
$ws.Cells.Item(3, "F").Value = 20.462
$ws.Cells.Item(3, "G").Value = 3.047
$ws.Cells.Item(3, "H").Value = 22.093
$ws.Cells.Item(3, "I").Value = 2.294
$ws.Cells.Item(3, "K").Value = 8.565
$ws.Cells.Item(3, "L").Value = 17.975
$ws.Cells.Item(3, "M").Value = 6.507

# Row 4
$ws.Cells.Item(4, "B").Value = 0
$ws.Cells.Item(4, "C").Value = 0
$ws.Cells.Item(4, "D").Value = 0
$ws.Cells.Item(4, "E").Value = 0
$ws.Cells.Item(4, "F").Value = 0
$ws.Cells.Item(4, "G").Value = 0
$ws.Cells.Item(4, "H").Value = 0
$ws.Cells.Item(4, "I").Value = 0
$ws.Cells.Item(4, "K").Value = 0
$ws.Cells.Item(4, "L").Value = 0
$ws.Cells.Item(4, "M").Value = 0

# Row 5
$ws.Cells.Item(5, "B").Value = -11.019
$ws.Cells.Item(5, "C").Value = 23.664
$ws.Cells.Item(5, "D").Value = -28.644
$ws.Cells.Item(5, "E").Value = -1.191
$ws.Cells.Item(5, "F").Value = 15.26
$ws.Cells.Item(5, "G").Value = 10.569
$ws.Cells.Item(5, "H").Value = -7.81
$ws.Cells.Item(5, "I").Value = 5.096
$ws.Cells.Item(5, "K").Value = 40.579
$ws.Cells.Item(5, "L").Value = 26.442
$ws.Cells.Item(5, "M").Value = 3.12

# Row 8
$ws.Cells.Item(8, "B").Value = 0.779
$ws.Cells.Item(8, "C").Value = 1.002
$ws.Cells.Item(8, "D").Value = 1.128
$ws.Cells.Item(8, "E").Value = 0.928
$ws.Cells.Item(8, "F").Value = 1.068
$ws.Cells.Item(8, "G").Value = 0.996
$ws.Cells.Item(8, "H").Value = 1.036
$ws.Cells.Item(8, "I").Value = 0.792
$ws.Cells.Item(8, "K").Value = 1.002
$ws.Cells.Item(8, "L").Value = 0.95
$ws.Cells.Item(8, "M").Value = 0.558

# Row 9
$ws.Cells.Item(9, "B").Value = -6.433
$ws.Cells.Item(9, "C").Value = 12.204
$ws.Cells.Item(9, "D").Value = -16.972
$ws.Cells.Item(9, "E").Value = 12.032
$ws.Cells.Item(9, "F").Value = 10.055
$ws.Cells.Item(9, "G").Value = 0.364
$ws.Cells.Item(9, "H").Value = -6.676
$ws.Cells.Item(9, "I").Value = 4.701
$ws.Cells.Item(9, "K").Value = 19.095
$ws.Cells.Item(9, "L").Value = 8.887
$ws.Cells.Item(9, "M").Value = -2.806

# Row 10
$ws.Cells.Item(10, "B").Value = -0.053
$ws.Cells.Item(10, "C").Value = 0.905
$ws.Cells.Item(10, "D").Value = -0.298
$ws.Cells.Item(10, "E").Value = 0.761
$ws.Cells.Item(10, "F").Value = 0.761
$ws.Cells.Item(10, "G").Value = 0.303
$ws.Cells.Item(10, "H").Value = 0.022
$ws.Cells.Item(10, "I").Value = 0.449
$ws.Cells.Item(10, "K").Value = 0.942
$ws.Cells.Item(10, "L").Value = 0.734
$ws.Cells.Item(10, "M").Value = 0.108

# Row 11
$ws.Cells.Item(11, "B").Value = -0.013
$ws.Cells.Item(11, "D").Value = -0.081
$ws.Cells.Item(11, "G").Value = 0.073
$ws.Cells.Item(11, "H").Value = 0.005
$ws.Cells.Item(11, "I").Value = 0.129
$ws.Cells.Item(11, "K").Value = 0.26
$ws.Cells.Item(11, "M").Value = 0.019

# Row 14
$ws.Cells.Item(14, "B").Value = 0.398
$ws.Cells.Item(14, "C").Value = 0.727
$ws.Cells.Item(14, "D").Value = 0.466
$ws.Cells.Item(14, "E").Value = 0.467
$ws.Cells.Item(14, "F").Value = 0.836
$ws.Cells.Item(14, "G").Value = 0.519
$ws.Cells.Item(14, "H").Value = 0.575
$ws.Cells.Item(14, "I").Value = 0.381
$ws.Cells.Item(14, "K").Value = 0.633
$ws.Cells.Item(14, "L").Value = 0.745
$ws.Cells.Item(14, "M").Value = 0.329

# Row 15
$ws.Cells.Item(15, "B").Value = 6.464
$ws.Cells.Item(15, "C").Value = 8.016
$ws.Cells.Item(15, "D").Value = 8.888
$ws.Cells.Item(15, "E").Value = 7.5
$ws.Cells.Item(15, "F").Value = 8.472
$ws.Cells.Item(15, "G").Value = 7.971
$ws.Cells.Item(15, "H").Value = 8.254
$ws.Cells.Item(15, "I").Value = 6.555
$ws.Cells.Item(15, "K").Value = 8.012
$ws.Cells.Item(15, "L").Value = 7.652
$ws.Cells.Item(15, "M").Value = 4.926
